$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column AA (27) mirrors column Z (26) for header rows 1-17
$headerValues = @(
    "age",
    "education",
    "income",
    "race",
    "sex",
    "release",
    "star_user",
    "real_extraversion",
    "real_agreeableness",
    "real_conscientiousness",
    "real_emotionstability",
    "real_openness",
    "gap_extraversion",
    "gap_agreeableness",
    "gap_conscientiousness",
    "gap_emotionstability",
    "gap_openness"
)

for ($i = 0; $i -lt $headerValues.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 27).Value = $headerValues[$i]
    $ws.Cells.Item($row, 28).Value = $headerValues[$i]
}

# Column AA (27): tste_17_* labels for rows 18-119
$tste17 = @(
    "tste_17_0",
    "tste_17_1",
    "tste_17_2",
    "tste_17_3",
    "tste_17_4",
    "tste_17_5",
    "tste_17_6",
    "tste_17_7",
    "tste_17_8",
    "tste_17_9",
    "tste_17_10",
    "tste_17_11",
    "tste_17_12",
    "tste_17_13",
    "tste_17_14",
    "tste_17_15",
    "tste_17_16",
    "tste_17_0*gap_extraversion",
    "tste_17_1*gap_extraversion",
    "tste_17_2*gap_extraversion",
    "tste_17_3*gap_extraversion",
    "tste_17_4*gap_extraversion",
    "tste_17_5*gap_extraversion",
    "tste_17_6*gap_extraversion",
    "tste_17_7*gap_extraversion",
    "tste_17_8*gap_extraversion",
    "tste_17_9*gap_extraversion",
    "tste_17_10*gap_extraversion",
    "tste_17_11*gap_extraversion",
    "tste_17_12*gap_extraversion",
    "tste_17_13*gap_extraversion",
    "tste_17_14*gap_extraversion",
    "tste_17_15*gap_extraversion",
    "tste_17_16*gap_extraversion",
    "tste_17_0*gap_agreeableness",
    "tste_17_1*gap_agreeableness",
    "tste_17_2*gap_agreeableness",
    "tste_17_3*gap_agreeableness",
    "tste_17_4*gap_agreeableness",
    "tste_17_5*gap_agreeableness",
    "tste_17_6*gap_agreeableness",
    "tste_17_7*gap_agreeableness",
    "tste_17_8*gap_agreeableness",
    "tste_17_9*gap_agreeableness",
    "tste_17_10*gap_agreeableness",
    "tste_17_11*gap_agreeableness",
    "tste_17_12*gap_agreeableness",
    "tste_17_13*gap_agreeableness",
    "tste_17_14*gap_agreeableness",
    "tste_17_15*gap_agreeableness",
    "tste_17_16*gap_agreeableness",
    "tste_17_0*gap_conscientiousness",
    "tste_17_1*gap_conscientiousness",
    "tste_17_2*gap_conscientiousness",
    "tste_17_3*gap_conscientiousness",
    "tste_17_4*gap_conscientiousness",
    "tste_17_5*gap_conscientiousness",
    "tste_17_6*gap_conscientiousness",
    "tste_17_7*gap_conscientiousness",
    "tste_17_8*gap_conscientiousness",
    "tste_17_9*gap_conscientiousness",
    "tste_17_10*gap_conscientiousness",
    "tste_17_11*gap_conscientiousness",
    "tste_17_12*gap_conscientiousness",
    "tste_17_13*gap_conscientiousness",
    "tste_17_14*gap_conscientiousness",
    "tste_17_15*gap_conscientiousness",
    "tste_17_16*gap_conscientiousness",
    "tste_17_0*gap_emotionstability",
    "tste_17_1*gap_emotionstability",
    "tste_17_2*gap_emotionstability",
    "tste_17_3*gap_emotionstability",
    "tste_17_4*gap_emotionstability",
    "tste_17_5*gap_emotionstability",
    "tste_17_6*gap_emotionstability",
    "tste_17_7*gap_emotionstability",
    "tste_17_8*gap_emotionstability",
    "tste_17_9*gap_emotionstability",
    "tste_17_10*gap_emotionstability",
    "tste_17_11*gap_emotionstability",
    "tste_17_12*gap_emotionstability",
    "tste_17_13*gap_emotionstability",
    "tste_17_14*gap_emotionstability",
    "tste_17_15*gap_emotionstability",
    "tste_17_16*gap_emotionstability",
    "tste_17_0*gap_openness",
    "tste_17_1*gap_openness",
    "tste_17_2*gap_openness",
    "tste_17_3*gap_openness",
    "tste_17_4*gap_openness",
    "tste_17_5*gap_openness",
    "tste_17_6*gap_openness",
    "tste_17_7*gap_openness",
    "tste_17_8*gap_openness",
    "tste_17_9*gap_openness",
    "tste_17_10*gap_openness",
    "tste_17_11*gap_openness",
    "tste_17_12*gap_openness",
    "tste_17_13*gap_openness",
    "tste_17_14*gap_openness",
    "tste_17_15*gap_openness",
    "tste_17_16*gap_openness"
)

for ($i = 0; $i -lt $tste17.Count; $i++) {
    $row = 18 + $i
    $ws.Cells.Item($row, 27).Value = $tste17[$i]
}

# Column AB (28): tste_18_* labels for rows 18-125
$tste18 = @(
    "tste_18_0",
    "tste_18_1",
    "tste_18_2",
    "tste_18_3",
    "tste_18_4",
    "tste_18_5",
    "tste_18_6",
    "tste_18_7",
    "tste_18_8",
    "tste_18_9",
    "tste_18_10",
    "tste_18_11",
    "tste_18_12",
    "tste_18_13",
    "tste_18_14",
    "tste_18_15",
    "tste_18_16",
    "tste_18_17",
    "tste_18_0*gap_extraversion",
    "tste_18_1*gap_extraversion",
    "tste_18_2*gap_extraversion",
    "tste_18_3*gap_extraversion",
    "tste_18_4*gap_extraversion",
    "tste_18_5*gap_extraversion",
    "tste_18_6*gap_extraversion",
    "tste_18_7*gap_extraversion",
    "tste_18_8*gap_extraversion",
    "tste_18_9*gap_extraversion",
    "tste_18_10*gap_extraversion",
    "tste_18_11*gap_extraversion",
    "tste_18_12*gap_extraversion",
    "tste_18_13*gap_extraversion",
    "tste_18_14*gap_extraversion",
    "tste_18_15*gap_extraversion",
    "tste_18_16*gap_extraversion",
    "tste_18_17*gap_extraversion",
    "tste_18_0*gap_agreeableness",
    "tste_18_1*gap_agreeableness",
    "tste_18_2*gap_agreeableness",
    "tste_18_3*gap_agreeableness",
    "tste_18_4*gap_agreeableness",
    "tste_18_5*gap_agreeableness",
    "tste_18_6*gap_agreeableness",
    "tste_18_7*gap_agreeableness",
    "tste_18_8*gap_agreeableness",
    "tste_18_9*gap_agreeableness",
    "tste_18_10*gap_agreeableness",
    "tste_18_11*gap_agreeableness",
    "tste_18_12*gap_agreeableness",
    "tste_18_13*gap_agreeableness",
    "tste_18_14*gap_agreeableness",
    "tste_18_15*gap_agreeableness",
    "tste_18_16*gap_agreeableness",
    "tste_18_17*gap_agreeableness",
    "tste_18_0*gap_conscientiousness",
    "tste_18_1*gap_conscientiousness",
    "tste_18_2*gap_conscientiousness",
    "tste_18_3*gap_conscientiousness",
    "tste_18_4*gap_conscientiousness",
    "tste_18_5*gap_conscientiousness",
    "tste_18_6*gap_conscientiousness",
    "tste_18_7*gap_conscientiousness",
    "tste_18_8*gap_conscientiousness",
    "tste_18_9*gap_conscientiousness",
    "tste_18_10*gap_conscientiousness",
    "tste_18_11*gap_conscientiousness",
    "tste_18_12*gap_conscientiousness",
    "tste_18_13*gap_conscientiousness",
    "tste_18_14*gap_conscientiousness",
    "tste_18_15*gap_conscientiousness",
    "tste_18_16*gap_conscientiousness",
    "tste_18_17*gap_conscientiousness",
    "tste_18_0*gap_emotionstability",
    "tste_18_1*gap_emotionstability",
    "tste_18_2*gap_emotionstability",
    "tste_18_3*gap_emotionstability",
    "tste_18_4*gap_emotionstability",
    "tste_18_5*gap_emotionstability",
    "tste_18_6*gap_emotionstability",
    "tste_18_7*gap_emotionstability",
    "tste_18_8*gap_emotionstability",
    "tste_18_9*gap_emotionstability",
    "tste_18_10*gap_emotionstability",
    "tste_18_11*gap_emotionstability",
    "tste_18_12*gap_emotionstability",
    "tste_18_13*gap_emotionstability",
    "tste_18_14*gap_emotionstability",
    "tste_18_15*gap_emotionstability",
    "tste_18_16*gap_emotionstability",
    "tste_18_17*gap_emotionstability",
    "tste_18_0*gap_openness",
    "tste_18_1*gap_openness",
    "tste_18_2*gap_openness",
    "tste_18_3*gap_openness",
    "tste_18_4*gap_openness",
    "tste_18_5*gap_openness",
    "tste_18_6*gap_openness",
    "tste_18_7*gap_openness",
    "tste_18_8*gap_openness",
    "tste_18_9*gap_openness",
    "tste_18_10*gap_openness",
    "tste_18_11*gap_openness",
    "tste_18_12*gap_openness",
    "tste_18_13*gap_openness",
    "tste_18_14*gap_openness",
    "tste_18_15*gap_openness",
    "tste_18_16*gap_openness",
    "tste_18_17*gap_openness"
)

for ($i = 0; $i -lt $tste18.Count; $i++) {
    $row = 18 + $i
    $ws.Cells.Item($row, 28).Value = $tste18[$i]
}

# Update sheet view: topLeftCell W1 -> actually diff says topLeftCell="W1", selection AC18
$ws.Application.ActiveWindow.ScrollColumn = 23
$ws.Range("AC18").Select()
